$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calibration data (concentration, peak area) pairs
$data = @(
    @(5,    0.0632),
    @(5,    0.0725),
    @(10,   0.1126),
    @(10,   0.1344),
    @(50,   0.6075),
    @(50,   0.583),
    @(100,  1.0714),
    @(100,  1.1227),
    @(500,  5.129),
    @(500,  5.4232),
    @(1000, 10.3892),
    @(1000, 10.5105),
    @(5000, 46.7262),
    @(5000, 51.1182)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("B4").Select()

# Update chart to point at the now-larger data range
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = $ws.Range("A2:A15")
$series.Values = $ws.Range("B2:B15")

